# Add proactive auto-enrollment and debug logging to AE process
# Update the projection results table with the new auto-enrollment figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ C=9401;  D=8311; E=0.884054887777896;  F=0.881710163377891;  G=0.1004018770304416;  H=0.08852535539995757; I=42616319.262861;  J=14844035.3774605; K=0; L=14844035.3774605; M=57460354.6403215; N=798582052.9413;   O=780882245.9373;   P=0.01858799020437242; Q=0.01900931344602805 }
    3 = @{ C=9592;  D=8469; E=0.8829232693911593; F=0.8809944866326849; G=0.09646357303105442; H=0.08498387600124831; I=45080697.01096167; J=16186053.48021406; K=0; L=16186053.48021406; M=61266750.49117575; N=844996599.4690611; O=827401163.456037;  P=0.01915516996208539; Q=0.01956252202088436 }
    4 = @{ C=9777;  D=8666; E=0.8863659609287102; F=0.8839249286005711; G=0.09398338333717977; H=0.08307425540595674; I=48078757.09076811; J=17612347.33775678; K=0; L=17612347.33775678; M=65691104.42852489; N=891229868.3884727; O=873665952.5938962; P=0.01976184591928403; Q=0.0201591320864297 }
    5 = @{ C=9972;  D=8833; E=0.8857801845166466; F=0.8833883388338833; G=0.09144571493263898; H=0.08078207820782078; I=50775057.45786561; J=18958537.25250198; K=0; L=18958537.25250198; M=69733594.71036758; N=936359078.2269156; O=918757855.9416611; P=0.02024708009282268; Q=0.02063496614466587 }
    6 = @{ C=10166; D=9018; E=0.8870745622663782; F=0.8842910374583252; G=0.08842980705256154; H=0.07819768582074915; I=53342426.2551169;  J=20351365.58954056; K=0; L=20351365.58954056; M=73693791.84465745; N=982369991.8607023; O=964662445.6605709; P=0.02071659940568129; Q=0.02109687764988568 }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
